$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected "new positive cases" input values (cumulative totals in column B
# are formula-driven and recalculate automatically from these).
$ws.Range("C231").Value = 222
$ws.Range("C241").Value = 760

# Corrected hospital/extra-hospital death split for a few existing rows.
$ws.Range("L263").Value = 3
$ws.Range("L264").Value = 4
$ws.Range("L267").Value = 4
$ws.Range("M267").Value = 8

$ws.Range("C268").Value = 195
$ws.Range("M268").Value = 3

$ws.Range("C269").Value = 150
$ws.Range("L269").Value = 4

# Newly-filled-in data rows 270-272 (previously blank placeholders).
$ws.Range("C270").Value = 117
$ws.Range("E270").Value = 34
$ws.Range("F270").Value = 24
$ws.Range("G270").Value = 199
$ws.Range("L270").Value = 5
$ws.Range("M270").Value = 1

$ws.Range("C271").Value = 73
$ws.Range("E271").Value = 34
$ws.Range("F271").Value = 24
$ws.Range("G271").Value = 194
$ws.Range("L271").Value = 1
$ws.Range("M271").Value = 0

$ws.Range("C272").Value = 29
$ws.Range("E272").Value = 34
$ws.Range("F272").Value = 26
$ws.Range("G272").Value = 201
$ws.Range("L272").Value = 0
$ws.Range("M272").Value = 0
